$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 11.7109375
$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 11.7109375

# Update cell values
$ws.Range("A1").Value = 149.01143520595124
$ws.Range("B1").Value = 4.7780736995762796
$ws.Range("C1").Value = 0.69070450097847369
